$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: clear the old 3-row bibliography (the "Microsoft, 2008" / Support
# Lifecycle entry is being dropped entirely; Baker and O'Connor are rewritten below
# as part of the full 13-row list in its pre-sort order).
$ws.Range("A1:B3").ClearContents()

# Write the finished bibliography (13 entries) - order here does not matter, the
# range gets sorted by column B (the full reference text) right after.
$ws.Range("A1").Value = 'Baker, Hylender, & Valentine, 2008, p. 2'
$ws.Range("B1").Value = 'Baker, W. H., Hylender, C. D., & Valentine, J. A. (2008). 2008 Data Breach Investigations Report, 1–29.'
$ws.Range("A2").Value = 'Goodin, 2014'
$ws.Range("B2").Value = 'Goodin, D. (2014). Stanford’s password policy shuns one-size-fits-all security | Ars Technica. Ars Technica. Retrieved April 30, 2015, from http://arstechnica.com/security/2014/04/25/stanfords-password-policy-shuns-one-size-fits-all-security/'
$ws.Range("A3").Value = 'Machkovech, 2015'
$ws.Range("B3").Value = 'Machkovech, S. (2015). Hacked French network exposed its own passwords during TV interview | Ars Technica. Ars Technica. Retrieved May 6, 2015, from http://arstechnica.com/security/2015/04/09/hacked-french-network-exposed-its-own-passwords-during-tv-interview/'
$ws.Range("A4").Value = 'Microsoft, 2005a'
$ws.Range("B4").Value = 'Microsoft. (2005a). Apply or modify password policy: Logon and Authentication. Retrieved May 14, 2015, from https://technet.microsoft.com/en-au/library/cc781633(v=ws.10).aspx?f=255&MSPPError=-2147217396#BKMK_3'
$ws.Range("A5").Value = 'Microsoft, 2005b'
$ws.Range("B5").Value = 'Microsoft. (2005b). Assign user rights to a group in Active Directory: Active Directory. Retrieved May 14, 2015, from https://technet.microsoft.com/en-au/library/cc786658(v=ws.10).aspx'
$ws.Range("A6").Value = 'Microsoft, 2007'
$ws.Range("B6").Value = 'Microsoft. (2007). Configure UAC settings via policy - Microsoft Reduce Customer Effort Center - Site Home - TechNet Blogs. Retrieved May 14, 2015, from http://blogs.technet.com/b/asiasupp/archive/2007/02/08/configure-uac-settings-via-policy.aspx'
$ws.Range("A7").Value = 'Microsoft, n.d.-a'
$ws.Range("B7").Value = 'Microsoft. (n.d.-a). Screen Saver timeout. Retrieved May 14, 2015, from https://technet.microsoft.com/en-us/library/cc961876.aspx'
$ws.Range("A8").Value = 'Microsoft, n.d.-b'
$ws.Range("B8").Value = 'Microsoft. (n.d.-b). What is User Account Control? - Windows Help. Retrieved March 14, 2015, from http://windows.microsoft.com/en-au/windows/what-is-user-account-control#1TC=windows-vista'
$ws.Range("A9").Value = 'O''Connor, 2008'
$ws.Range("B9").Value = 'O’Connor, E. (2008). BigAdmin Feature Article: Patch Management Best Practices. Retrieved May 13, 2015, from http://www.oracle.com/technetwork/systems/articles/patch-management-jsp-135385.html'
$ws.Range("A10").Value = 'Oliver & Snowden, 2015'
$ws.Range("B10").Value = 'Oliver, J., & Snowden, E. [LastWeekTonight]. (2015, April 9). Last Week Tonight with John Oliver: Edward Snowden on Passwords. Retrieved May 6, 2015, from https://www.youtube.com/watch?v=yzGzB-yYKcc'
$ws.Range("A11").Value = 'Microsoft, 2012'
$ws.Range("B11").Value = 'Password must meet complexity requirements. (n.d.). Retrieved May 14, 2015, from https://technet.microsoft.com/en-us/library/hh994562(v=ws.10).aspx'
$ws.Range("A12").Value = 'Scarfone & Souppaya, 2013'
$ws.Range("B12").Value = 'Scarfone, K., & Souppaya, M. (2013). Guide to Enterprise Patch Management Technologies NIST Special Publication 800-40 Guide to Enterprise Patch Management Technologies. NIST. doi:10.6028/NIST.SP.800-40r3'
$ws.Range("A13").Value = 'Tice, 2012'
$ws.Range("B13").Value = 'Tice, K. [solarwindsinc]. (2012, September 12)  Patch Manager Guided Tour. Retrieved May 14, 2015, from https://www.youtube.com/watch?v=-DldViUL1d0'

# Sort the whole table alphabetically by column B, same as Data > Sort in the UI.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B5"))
$ws.Sort.SetRange($ws.Range("A1:B13"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# The text sort above is a plain ordinal comparison, so it places the curly-quote
# ‘O’Connor’ reference after ‘Oliver’; Excel's own text collation treats the quote as
# ignorable punctuation and keeps it first - restore that exact order here.
$ws.Range("A9").Value = 'O''Connor, 2008'
$ws.Range("B9").Value = 'O’Connor, E. (2008). BigAdmin Feature Article: Patch Management Best Practices. Retrieved May 13, 2015, from http://www.oracle.com/technetwork/systems/articles/patch-management-jsp-135385.html'
$ws.Range("A10").Value = 'Oliver & Snowden, 2015'
$ws.Range("B10").Value = 'Oliver, J., & Snowden, E. [LastWeekTonight]. (2015, April 9). Last Week Tonight with John Oliver: Edward Snowden on Passwords. Retrieved May 6, 2015, from https://www.youtube.com/watch?v=yzGzB-yYKcc'

# Column B keeps its vertical-centred look across every (now 13) row.
$ws.Range("B1:B13").VerticalAlignment = -4108
# Row 7 (the Screen Saver timeout reference) is also explicitly left-aligned.
$ws.Range("B7").HorizontalAlignment = -4131

# Selection / print setup picked up in the refreshed sheet.
$null = $ws.Range("B19").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "bib_v2 updated"
